$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 272.8889
$ws.Range("I41").Value = 161.33333
$ws.Range("J41").Value = 384.44446
$ws.Range("K41").Value = 161.33333
$ws.Range("L41").Value = 384.44446
$ws.Range("M41").Value = 278.66667
$ws.Range("N41").Value = -1264.44446

$ws.Range("H92").Value = 562.93335
$ws.Range("I92").Value = 335.7143
$ws.Range("J92").Value = 761.75
$ws.Range("K92").Value = 335.7143
$ws.Range("L92").Value = 761.75
$ws.Range("M92").Value = 912.2857
$ws.Range("N92").Value = -3257.75

$ws.Range("H101").Value = 276.85715
$ws.Range("I101").Value = 241.33333
$ws.Range("J101").Value = 490
$ws.Range("K101").Value = 723.99999
$ws.Range("L101").Value = 1470
$ws.Range("M101").Value = 898.00001
$ws.Range("N101").Value = -4714

$ws.Range("H131").Value = 78118.46000000001
$ws.Range("I131").Value = 84336.664
$ws.Range("K131").Value = 253009.992
$ws.Range("M131").Value = -247969.992

$ws.Range("H141").Value = 5129.0625
$ws.Range("I141").Value = 3247.7778
$ws.Range("K141").Value = 9743.3334
$ws.Range("M141").Value = -4563.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 118.35294
$ws.Range("I5").Value = 126.72727
$ws.Range("K5").Value = 126.72727
$ws.Range("M5").Value = -14.72727

$ws.Range("H45").Value = 1100
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -1954

$ws.Range("H121").Value = 23993.334
$ws.Range("J121").Value = 23993.334
$ws.Range("L121").Value = 23993.334
$ws.Range("N121").Value = -27487.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 118.35294
$ws.Range("I4").Value = 126.72727
$ws.Range("K4").Value = 126.72727
$ws.Range("M4").Value = -11.72727

$ws.Range("H15").Value = 29666.666
$ws.Range("I15").Value = 9000
$ws.Range("J15").Value = 40000
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 40000
$ws.Range("M15").Value = -8773
$ws.Range("N15").Value = -40454

$ws.Range("H82").Value = 24639.625
$ws.Range("I82").Value = 8557
$ws.Range("J82").Value = 26937.143
$ws.Range("K82").Value = 8557
$ws.Range("L82").Value = 26937.143
$ws.Range("M82").Value = -8174
$ws.Range("N82").Value = -27703.143

$ws.Range("H85").Value = 24639.625
$ws.Range("I85").Value = 8557
$ws.Range("J85").Value = 26937.143
$ws.Range("K85").Value = 8557
$ws.Range("L85").Value = 26937.143
$ws.Range("M85").Value = -7231
$ws.Range("N85").Value = -29589.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 400
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 100
$ws.Range("M6").Value = -887
$ws.Range("N6").Value = -326

$ws.Range("H17").Value = 28998
$ws.Range("I17").Value = 28998
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 28998
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -28824
$ws.Range("N17").ClearContents()

$ws.Range("H25").Value = 15375
$ws.Range("I25").Value = 8833.333000000001
$ws.Range("K25").Value = 8833.333000000001
$ws.Range("M25").Value = -8659.333000000001

$ws.Range("H31").Value = 9806299
$ws.Range("I31").Value = 13515020
$ws.Range("J31").Value = 4679.2144
$ws.Range("K31").Value = 13515020
$ws.Range("L31").Value = 4679.2144
$ws.Range("M31").Value = -13514725
$ws.Range("N31").Value = -5269.2144

$ws.Range("H34").Value = 9806299
$ws.Range("I34").Value = 13515020
$ws.Range("J34").Value = 4679.2144
$ws.Range("K34").Value = 13515020
$ws.Range("L34").Value = 4679.2144
$ws.Range("M34").Value = -13514818
$ws.Range("N34").Value = -5083.2144

$ws.Range("H41").Value = 31156
$ws.Range("I41").Value = 7000
$ws.Range("J41").Value = 37195
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 37195
$ws.Range("M41").Value = -6572
$ws.Range("N41").Value = -38051

$ws.Range("H50").Value = 11281
$ws.Range("J50").Value = 11281
$ws.Range("L50").Value = 11281
$ws.Range("N50").Value = -12531

$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11472
$ws.Range("M51").ClearContents()

$ws.Range("H59").Value = 39956
$ws.Range("J59").Value = 39956
$ws.Range("L59").Value = 39956
$ws.Range("N59").Value = -42246

$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10696
$ws.Range("M61").ClearContents()

$ws.Range("H68").Value = 39525.215
$ws.Range("I68").Value = 27268
$ws.Range("J68").Value = 40468.08
$ws.Range("K68").Value = 27268
$ws.Range("L68").Value = 40468.08
$ws.Range("M68").Value = -26519
$ws.Range("N68").Value = -41966.08

$ws.Range("H71").Value = 39525.215
$ws.Range("I71").Value = 27268
$ws.Range("J71").Value = 40468.08
$ws.Range("K71").Value = 81804
$ws.Range("L71").Value = 121404.24
$ws.Range("M71").Value = -78060
$ws.Range("N71").Value = -128892.24

$ws.Range("H74").Value = 14171.272
$ws.Range("J74").Value = 14171.272
$ws.Range("L74").Value = 14171.272
$ws.Range("N74").Value = -15919.272

$ws.Range("H77").Value = 14171.272
$ws.Range("J77").Value = 14171.272
$ws.Range("L77").Value = 42513.81600000001
$ws.Range("N77").Value = -51249.81600000001

$ws.Range("H99").Value = 2254.2856
$ws.Range("I99").Value = 1720
$ws.Range("J99").Value = 2521.4285
$ws.Range("K99").Value = 1720
$ws.Range("L99").Value = 2521.4285
$ws.Range("M99").Value = -222
$ws.Range("N99").Value = -5517.4285

$ws.Range("H126").Value = 2254.2856
$ws.Range("I126").Value = 1720
$ws.Range("J126").Value = 2521.4285
$ws.Range("K126").Value = 5160
$ws.Range("L126").Value = 7564.2855
$ws.Range("M126").Value = -2690
$ws.Range("N126").Value = -12504.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 13158737
$ws.Range("I121").Value = 214.5
$ws.Range("J121").Value = 14706799
$ws.Range("K121").Value = 643.5
$ws.Range("L121").Value = 44120397
$ws.Range("M121").Value = 666.5
$ws.Range("N121").Value = -44123017

$ws.Range("H122").Value = 1014.36365
$ws.Range("I122").Value = 876.5
$ws.Range("J122").Value = 1179.8
$ws.Range("K122").Value = 7888.5
$ws.Range("L122").Value = 10618.2
$ws.Range("M122").Value = -5438.5
$ws.Range("N122").Value = -15518.2

$ws.Range("H131").Value = 925.6923
$ws.Range("J131").Value = 969.5
$ws.Range("L131").Value = 2908.5
$ws.Range("N131").Value = -12988.5

$ws.Range("H139").Value = 4105.1904
$ws.Range("I139").Value = 4145.4443
$ws.Range("J139").Value = 4075
$ws.Range("K139").Value = 12436.3329
$ws.Range("L139").Value = 12225
$ws.Range("M139").Value = -7296.332900000001
$ws.Range("N139").Value = -22505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 4543.1665
$ws.Range("I132").Value = 4768.567
$ws.Range("J132").Value = 3416.1667
$ws.Range("K132").Value = 14305.701
$ws.Range("L132").Value = 10248.5001
$ws.Range("M132").Value = -11775.701
$ws.Range("N132").Value = -15308.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3833.3333
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3833.3333
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3833.3333
$ws.Range("N100").Value = -4915.3333
$ws.Range("M100").ClearContents()

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 5125.2
$ws.Range("I122").Value = 5353.1763
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 16059.5289
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -13609.5289
$ws.Range("N122").Value = -16399.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 33600
$ws.Range("J64").Value = 33600
$ws.Range("L64").Value = 33600
$ws.Range("N64").Value = -34096

$ws.Range("H67").Value = 33600
$ws.Range("J67").Value = 33600
$ws.Range("L67").Value = 33600
$ws.Range("N67").Value = -35316

$ws.Range("H82").Value = 35301
$ws.Range("J82").Value = 35301
$ws.Range("L82").Value = 35301
$ws.Range("N82").Value = -36067

$ws.Range("H85").Value = 35301
$ws.Range("J85").Value = 35301
$ws.Range("L85").Value = 35301
$ws.Range("N85").Value = -37953
